$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "device_master" rows for two additional MAC-address batches (set 30 and set 31) ---
# Column layout: A=id, B=name, C=mac_address, D=serial_num, E=(blank ip_address),
#                F=dspec_id, G=lang_code, H=is_active, I=cr_by, J=cr_dtimes

$idCol     = 1
$nameCol   = 2
$macCol    = 3
$serialCol = 4
$dspecCol  = 6
$langCol   = 7
$activeCol = 8
$crByCol   = 9
$crDtCol   = 10

$batch30 = @(
    @{ Row = 147; Id = 3000166; Name = "Finger Print Scanner 30"; Mac = "D6-15-AC-80-6B-86"; Serial = "BS563Q2230814"; Dspec = 165 },
    @{ Row = 148; Id = 3000167; Name = "IRIS Scanner 30";         Mac = "6D-58-E2-DF-74-34"; Serial = "BS563Q2230815"; Dspec = 327 },
    @{ Row = 149; Id = 3000168; Name = "Web Camera 30";           Mac = "E2-A8-56-86-15-30"; Serial = "BS563Q2230816"; Dspec = 736 },
    @{ Row = 150; Id = 3000169; Name = "Document Scanner 30";     Mac = "72-E8-B9-FD-63-65"; Serial = "BS563Q2230817"; Dspec = 801 },
    @{ Row = 151; Id = 3000170; Name = "Printer 30";              Mac = "D3-F3-A4-50-AD-12"; Serial = "BS563Q2230818"; Dspec = 920 }
)

$batch31 = @(
    @{ Row = 152; Id = 3000171; Name = "Finger Print Scanner 31"; Mac = "06-16-D0-0B-A6-E4"; Serial = "BS563Q2230819"; Dspec = 165 },
    @{ Row = 153; Id = 3000172; Name = "IRIS Scanner 31";         Mac = "21-78-45-AC-E9-20"; Serial = "BS563Q2230820"; Dspec = 327 },
    @{ Row = 154; Id = 3000173; Name = "Web Camera 31";           Mac = "3C-E8-87-99-DB-FA"; Serial = "BS563Q2230821"; Dspec = 736 },
    @{ Row = 155; Id = 3000174; Name = "Document Scanner 31";     Mac = "BF-55-53-98-40-08"; Serial = "BS563Q2230822"; Dspec = 801 },
    @{ Row = 156; Id = 3000175; Name = "Printer 31";              Mac = "5A-43-36-46-22-EB"; Serial = "BS563Q2230823"; Dspec = 920 }
)

# Fill column-by-column (name, then mac, then serial) per batch so new shared-string
# entries land in the same order the source workbook recorded them in.
foreach ($r in $batch30) { $ws.Cells.Item($r.Row, $nameCol).Value = $r.Name }
foreach ($r in $batch30) { $ws.Cells.Item($r.Row, $macCol).Value = $r.Mac }
foreach ($r in $batch30) { $ws.Cells.Item($r.Row, $serialCol).Value = $r.Serial }

foreach ($r in $batch31) { $ws.Cells.Item($r.Row, $nameCol).Value = $r.Name }
foreach ($r in $batch31) { $ws.Cells.Item($r.Row, $serialCol).Value = $r.Serial }
foreach ($r in $batch31) { $ws.Cells.Item($r.Row, $macCol).Value = $r.Mac }

foreach ($r in ($batch30 + $batch31)) {
    $ws.Cells.Item($r.Row, $idCol).Value = $r.Id
    $ws.Cells.Item($r.Row, $dspecCol).Value = $r.Dspec
    $ws.Cells.Item($r.Row, $langCol).Value = "eng"
    $ws.Cells.Item($r.Row, $activeCol).Value = $true
    $ws.Cells.Item($r.Row, $activeCol).HorizontalAlignment = -4131
    $ws.Cells.Item($r.Row, $crByCol).Value = "superadmin"
    $ws.Cells.Item($r.Row, $crDtCol).Value = "now()"
}

# --- Update the view so the newly-added tail of the sheet is in frame ---
$excel.ActiveWindow.ScrollRow = 139
$ws.Range("E156").Select()
